# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 5 (pushing the existing rows 5-42
# down to 6-43) and populate it with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 5:42 down to 6:43, duplicating row 5's formatting for the
# freshly inserted row (matches Excel's default "insert" behaviour).
$ws.Rows("5:5").Insert()

# Populate the new row 5 with the new weekly record.
$ws.Range("A5").Value2 = 11
$ws.Range("B5").Value2 = "Vega Monumental Concepción"
$ws.Range("C5").Value2 = "Bíobío"
$ws.Range("D5").Value2 = 44624
$ws.Range("E5").Value2 = 8
$ws.Range("F5").Value2 = 100112031
$ws.Range("G5").Value2 = "Poroto verde"
$ws.Range("H5").Value2 = "Magnum"
$ws.Range("I5").Value2 = "Primera"
$ws.Range("J5").Value2 = 110
$ws.Range("K5").Value2 = 26000
$ws.Range("L5").Value2 = 27000
$ws.Range("M5").Value2 = 26545
$ws.Range("N5").Value2 = "$/saco 25 kilos"
$ws.Range("O5").Value2 = "Región Metropolitana"
$ws.Range("P5").Value2 = 1062
$ws.Range("Q5").Value2 = 25
$ws.Range("R5").Value2 = "Hortaliza"
